$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values on row 2 (PC_gestionDocumental environment / test data) ---

# Ambiente / URL / Contrasenia / NroCuenta -> new environment + credentials
# A2 keeps its quotePrefix style (leading apostrophe forces quotePrefix, matching original formatting)
$ws.Range("A2").Value = "'i-preproducciongestion.segurossura.com.ar"
$ws.Range("B2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("D2").Value = "silverarrow"

# NroCuenta
$ws.Range("E2").Value = 5934358994

# Anio
$ws.Range("R2").Value = 2021

# Patente / Motor / Chasis
$ws.Range("W2").Value = "RGA011"
$ws.Range("X2").Value = "1234567RGA010"
$ws.Range("Y2").Value = "1234567RGA010"

# --- Remove the hyperlink that was attached to B2 (URL cell) ---
$ws.Range("B2").Hyperlinks.Delete() | Out-Null

# --- Update the sheet view / selection (scroll back to A1, select F8) ---
$ws.Range("F8").Select() | Out-Null
